$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Append three new monthly rows (81-83, for 2021-08-01 / 09-01 / 10-01)
#    to every data sheet, each keeping the sheet's existing indicator title
#    in column A and the Germany/Finland/Sweden figures in C/D/E.
# ---------------------------------------------------------------------------

$newRowsBySheet = @{
    1 = @(
        @(81, 44409, 5.5, 2.3, 9),
        @(82, 44440, 8.6, 4.5, 10.1),
        @(83, 44470, 10.2, 8.699999999999999, 16.4)
    )
    2 = @(
        @(81, 44409, -3.5, 4.2, 7.3),
        @(82, 44440, -1.1, 6.3, 6.6),
        @(83, 44470, -1.5, 4.5, 5.7)
    )
    3 = @(
        @(81, 44409, 23.9, 23.4, 30.7),
        @(82, 44440, 25.6, 22.3, 28.5),
        @(83, 44470, 23, 24, 32.7)
    )
    4 = @(
        @(81, 44409, 0.3, 7.2, 20),
        @(82, 44440, 0.7, 5.3, 23.5),
        @(83, 44470, -0.2, 10.2, 24.9)
    )
    5 = @(
        @(81, 44409, 21, 13.5, 38.1),
        @(82, 44440, 18, 15.9, 38.7),
        @(83, 44470, 19.8, 19.3, 42.6)
    )
    6 = @(
        @(81, 44409, 117.2, 116.1, 123.2),
        @(82, 44440, 118, 116.8, 121.7),
        @(83, 44470, 117.5, 117, 123.8)
    )
    7 = @(
        @(81, 44409, 113.5, 111, 114.5),
        @(82, 44440, 114.4, 116, 115.7),
        @(83, 44470, 112.8, 117.7, 119.7)
    )
}

# Each sheet repeats a single fixed indicator name down the whole column A.
$titleBySheet = @{
    1 = "Construction confidence indicator (5%)"
    2 = "Consumer confidence indicator (20%)"
    3 = "Industrial confidence indicator (40%)"
    4 = "Retail trade confidence indicator (5%)"
    5 = "Services confidence indicator (30 %)"
    6 = "The Economic sentiment indicator is a composite measure (average = 100)"
    7 = "The Employment expectations indicator is a composite measure (average = 100)"
}

foreach ($sheetIndex in 1..7) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $title = $titleBySheet[$sheetIndex]

    $rows = $newRowsBySheet[$sheetIndex]
    foreach ($row in $rows) {
        $r = $row[0]
        $ws.Cells.Item($r, 1).Value = $title
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 2).NumberFormat = "yyyy-mm-dd"
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
    }
}

# ---------------------------------------------------------------------------
# 2) Small revisions to previously published figures.
# ---------------------------------------------------------------------------

# Sheet 4 ("Retail trade confidence indicator (5%)"): column C (Germany)
# revisions for several historical rows.
$ws4 = $wb.Worksheets.Item(4)
$sheet4Corrections = @{
    2  = -9.5
    3  = -9
    4  = -7.6
    5  = -7
    6  = -5.7
    7  = -5
    8  = -1.4
    9  = 1.7
    22 = -3.7
    23 = -3.6
    24 = -2.3
    25 = -3.5
    26 = -5.5
    29 = -0.2
    31 = -1.6
    32 = -2.9
}
foreach ($r in $sheet4Corrections.Keys) {
    $ws4.Cells.Item($r, 3).Value = $sheet4Corrections[$r]
}

# Sheet 6 ("The Economic sentiment indicator ...") row 9, column C (Germany).
$ws6 = $wb.Worksheets.Item(6)
$ws6.Cells.Item(9, 3).Value = 106.5
